$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 11.64

$ws.Range("D3").Value = 10.36
$ws.Range("E3").Value = 10.31

$ws.Range("B4").Value = 8.359999999999999
$ws.Range("C4").Value = 9.640000000000001
$ws.Range("E4").Value = 10.01

$ws.Range("C5").Value = 9.69
$ws.Range("D5").Value = 9.99
$ws.Range("F5").Value = 10.06
$ws.Range("G5").Value = 9.49

$ws.Range("E6").Value = 9.94
$ws.Range("I6").Value = 10.7
$ws.Range("J6").Value = 6.5

$ws.Range("E7").Value = 10.51

$ws.Range("I8").Value = 8.17

$ws.Range("F9").Value = 9.300000000000001
$ws.Range("H9").Value = 11.83

$ws.Range("F10").Value = 13.5
